$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3407354241124098
$ws.Range("C2").Value = 0.5844350863438295
$ws.Range("D2").Value = 0.4442199068852112
$ws.Range("E2").Value = 0.6664982422221466
$ws.Range("F2").Value = 0.5885134887136687
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.1726228500093137
$ws.Range("C3").Value = 0.4551251828485784
$ws.Range("D3").Value = 0.3113871153733175
$ws.Range("E3").Value = 0.5580207123156966
$ws.Range("F3").Value = 0.5460334188492619
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.2034452847961829
$ws.Range("C4").Value = 0.3989700209739589
$ws.Range("D4").Value = 0.2615574828328119
$ws.Range("E4").Value = 0.5114269085928232
$ws.Range("F4").Value = 0.483661004847352
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.3635533387143693
$ws.Range("C5").Value = 0.3920234804570175
$ws.Range("D5").Value = 0.2041752638615354
$ws.Range("E5").Value = 0.4518575703266854
$ws.Range("F5").Value = 0.2771362769876616
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.3372887961726861
$ws.Range("C6").Value = 0.367293011109447
$ws.Range("D6").Value = 0.1918266328756511
$ws.Range("E6").Value = 0.4379801740668761
$ws.Range("F6").Value = 0.289203773041539
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.3465491651461942
$ws.Range("C7").Value = 0.3849987877080897
$ws.Range("D7").Value = 0.2040158891459988
$ws.Range("E7").Value = 0.4516811808632266
$ws.Range("F7").Value = 0.3006242113637049
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.3728115640212805
$ws.Range("C8").Value = 0.4172293592727421
$ws.Range("D8").Value = 0.2174429248205509
$ws.Range("E8").Value = 0.4663077576242443
$ws.Range("F8").Value = 0.2915344481050866
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.4140746927647911
$ws.Range("C9").Value = 0.4327887006433242
$ws.Range("D9").Value = 0.2237006186653438
$ws.Range("E9").Value = 0.4729699976376343
$ws.Range("F9").Value = 0.2387302033153843
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.3853625486457603
$ws.Range("C10").Value = 0.4009980026591189
$ws.Range("D10").Value = 0.1935302347147137
$ws.Range("E10").Value = 0.4399207141232539
$ws.Range("F10").Value = 0.2225500727871226
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.3626069129865366
$ws.Range("C11").Value = 0.3740635702753569
$ws.Range("D11").Value = 0.1753426339092577
$ws.Range("E11").Value = 0.4187393388604153
$ws.Range("F11").Value = 0.2207534083381827
$ws.Range("G11").Value = 10
